$d = $word.ActiveDocument

# Directeur technique : 8160 -> 10200
$d.Content.Find.Execute("8160", $true, $true, $false, $false, $false,
                         $true, 1, $false, "10200", 2)

# Chef de projet : 6100 -> 4270
$d.Content.Find.Execute("6100", $true, $true, $false, $false, $false,
                         $true, 1, $false, "4270", 2)

# Montant total estimé : 14260 -> 14470
$d.Content.Find.Execute("14260", $true, $true, $false, $false, $false,
                         $true, 1, $false, "14470", 2)

# Montant total : 56805 -> 57015
$d.Content.Find.Execute("56805", $true, $true, $false, $false, $false,
                         $true, 1, $false, "57015", 2)
